$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.398.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.570.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.493"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.796.77"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.578.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.400.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.514"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "

$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  +3.86%  "

$ws.Range("E32").Value = "  -2.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.378.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "

$ws.Range("E36").Value = "  +4.58%  "

$ws.Range("E37").Value = "  -2.01%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("E39").Value = "  +1.13%  "

$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("E41").Value = "  -1.68%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("E43").Value = "  +2.62%  "

$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0474"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "

$ws.Range("E46").Value = "  -4.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("E48").Value = "  -6.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.708.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "85.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.68%  "
